$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header column in H1, reusing the formatting of the other
# header cells (bold, centered, thin border - style index 1 in this sheet).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the corresponding data value for row 2
$ws.Range("H2").Value = 1
